$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(54, 8).Value = 4239.8  # H54: 4700 -> 4239.8
$ws.Cells.Item(54, 9).Value = 1199  # I54: 3800 -> 1199
$ws.Cells.Item(54, 11).Value = 1199  # K54: 3800 -> 1199
$ws.Cells.Item(54, 13).Value = -713  # M54: -3314 -> -713

$ws.Cells.Item(64, 8).Value = 4213  # H64: 4297.3335 -> 4213
$ws.Cells.Item(64, 10).Value = 3980  # J64: 4000 -> 3980
$ws.Cells.Item(64, 12).Value = 3980  # L64: 4000 -> 3980
$ws.Cells.Item(64, 14).Value = -4476  # N64: -4496 -> -4476

$ws.Cells.Item(67, 8).Value = 4213  # H67: 4297.3335 -> 4213
$ws.Cells.Item(67, 10).Value = 3980  # J67: 4000 -> 3980
$ws.Cells.Item(67, 12).Value = 3980  # L67: 4000 -> 3980
$ws.Cells.Item(67, 14).Value = -5696  # N67: -5716 -> -5696

$ws.Cells.Item(86, 8).Value = 2222  # H86: 3529.75 -> 2222
$ws.Cells.Item(86, 9).Value = 2570  # I86: 4593.2856 -> 2570
$ws.Cells.Item(86, 10).Value = 1468  # J86: 2040.8 -> 1468
$ws.Cells.Item(86, 11).Value = 2570  # K86: 4593.2856 -> 2570
$ws.Cells.Item(86, 12).Value = 1468  # L86: 2040.8 -> 1468
$ws.Cells.Item(86, 13).Value = -1447  # M86: -3470.2856 -> -1447
$ws.Cells.Item(86, 14).Value = -3714  # N86: -4286.8 -> -3714

$ws.Cells.Item(89, 8).Value = 2222  # H89: 3529.75 -> 2222
$ws.Cells.Item(89, 9).Value = 2570  # I89: 4593.2856 -> 2570
$ws.Cells.Item(89, 10).Value = 1468  # J89: 2040.8 -> 1468
$ws.Cells.Item(89, 11).Value = 12850  # K89: 22966.428 -> 12850
$ws.Cells.Item(89, 12).Value = 7340  # L89: 10204 -> 7340
$ws.Cells.Item(89, 13).Value = -7234  # M89: -17350.428 -> -7234
$ws.Cells.Item(89, 14).Value = -18572  # N89: -21436 -> -18572

$ws.Cells.Item(112, 8).Value = 2007.9445  # H112: 2123.5625 -> 2007.9445
$ws.Cells.Item(112, 10).Value = 2196.3333  # J112: 2367.6155 -> 2196.3333
$ws.Cells.Item(112, 12).Value = 6588.999899999999  # L112: 7102.8465 -> 6588.999899999999
$ws.Cells.Item(112, 14).Value = -8804.999899999999  # N112: -9318.8465 -> -8804.999899999999

$ws.Cells.Item(134, 8).Value = 37745  # H134: 37797.273 -> 37745
$ws.Cells.Item(134, 10).Value = 37745  # J134: 37797.273 -> 37745
$ws.Cells.Item(134, 12).Value = 37745  # L134: 37797.273 -> 37745
$ws.Cells.Item(134, 14).Value = -47885  # N134: -47937.273 -> -47885

$ws.Cells.Item(135, 8).Value = 38463340  # H135: 37038804 -> 38463340
$ws.Cells.Item(135, 9).Value = 712.73334  # I135: 722.5625 -> 712.73334
$ws.Cells.Item(135, 11).Value = 6414.60006  # K135: 6503.0625 -> 6414.60006
$ws.Cells.Item(135, 13).Value = -3879.60006  # M135: -3968.0625 -> -3879.60006

$ws.Cells.Item(136, 8).Value = 42006.152  # H136: 41461.332 -> 42006.152
$ws.Cells.Item(136, 10).Value = 42006.152  # J136: 41461.332 -> 42006.152
$ws.Cells.Item(136, 12).Value = 42006.152  # L136: 41461.332 -> 42006.152
$ws.Cells.Item(136, 14).Value = -52206.152  # N136: -51661.332 -> -52206.152

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(9, 8).Value = 10000  # H9: 9999.666999999999 -> 10000
$ws.Cells.Item(9, 10).Value = 10000  # J9: 9999.666999999999 -> 10000
$ws.Cells.Item(9, 12).Value = 10000  # L9: 9999.666999999999 -> 10000
$ws.Cells.Item(9, 14).Value = -10340  # N9: -10339.667 -> -10340

$ws.Cells.Item(20, 8).Value = 10000  # H20: 9999.666999999999 -> 10000
$ws.Cells.Item(20, 10).Value = 10000  # J20: 9999.666999999999 -> 10000
$ws.Cells.Item(20, 12).Value = 10000  # L20: 9999.666999999999 -> 10000
$ws.Cells.Item(20, 14).Value = -10540  # N20: -10539.667 -> -10540

$ws.Cells.Item(32, 8).Value = 7800.98  # H32: 9803.020500000001 -> 7800.98
$ws.Cells.Item(32, 9).Value = 6214.107  # I32: 7142.325 -> 6214.107
$ws.Cells.Item(32, 10).Value = 16132.0625  # J32: 21005.947 -> 16132.0625
$ws.Cells.Item(32, 11).Value = 6214.107  # K32: 7142.325 -> 6214.107
$ws.Cells.Item(32, 12).Value = 16132.0625  # L32: 21005.947 -> 16132.0625
$ws.Cells.Item(32, 13).Value = -5927.107  # M32: -6855.325 -> -5927.107
$ws.Cells.Item(32, 14).Value = -16706.0625  # N32: -21579.947 -> -16706.0625

$ws.Cells.Item(61, 8).Value = 100001550  # H61: 83334780 -> 100001550
$ws.Cells.Item(61, 9).Value = 166667500  # I61: 125000860 -> 166667500
$ws.Cells.Item(61, 11).Value = 166667500  # K61: 125000860 -> 166667500
$ws.Cells.Item(61, 13).Value = -166667288  # M61: -125000648 -> -166667288

$ws.Cells.Item(74, 8).Value = 1877.9688  # H74: 1844.9395 -> 1877.9688
$ws.Cells.Item(74, 9).Value = 1053.8182  # I74: 1042.2609 -> 1053.8182
$ws.Cells.Item(74, 11).Value = 1053.8182  # K74: 1042.2609 -> 1053.8182
$ws.Cells.Item(74, 13).Value = -179.8181999999999  # M74: -168.2609 -> -179.8181999999999

$ws.Cells.Item(77, 8).Value = 1877.9688  # H77: 1844.9395 -> 1877.9688
$ws.Cells.Item(77, 9).Value = 1053.8182  # I77: 1042.2609 -> 1053.8182
$ws.Cells.Item(77, 11).Value = 5269.090999999999  # K77: 5211.3045 -> 5269.090999999999
$ws.Cells.Item(77, 13).Value = -901.0909999999994  # M77: -843.3045000000002 -> -901.0909999999994

$ws.Cells.Item(88, 8).Value = 1425.5  # H88: 2000 -> 1425.5
$ws.Cells.Item(88, 10).Value = 1425.5  # J88: 2000 -> 1425.5
$ws.Cells.Item(88, 12).Value = 1425.5  # L88: 2000 -> 1425.5
$ws.Cells.Item(88, 14).Value = -2237.5  # N88: -2812 -> -2237.5

$ws.Cells.Item(91, 8).Value = 1425.5  # H91: 2000 -> 1425.5
$ws.Cells.Item(91, 10).Value = 1425.5  # J91: 2000 -> 1425.5
$ws.Cells.Item(91, 12).Value = 1425.5  # L91: 2000 -> 1425.5
$ws.Cells.Item(91, 14).Value = -4233.5  # N91: -4808 -> -4233.5

$ws.Cells.Item(132, 8).Value = 2655.98  # H132: 2702.0408 -> 2655.98
$ws.Cells.Item(132, 9).Value = 2145.8386  # I132: 2204.0667 -> 2145.8386
$ws.Cells.Item(132, 11).Value = 6437.5158  # K132: 6612.2001 -> 6437.5158
$ws.Cells.Item(132, 13).Value = -3907.5158  # M132: -4082.2001 -> -3907.5158

$ws.Cells.Item(135, 8).Value = 23333  # H135: 0 -> 23333
$ws.Cells.Item(135, 10).Value = 23333  # J135: 0 -> 23333
$ws.Cells.Item(135, 12).Value = 23333  # L135: 0 -> 23333
$ws.Cells.Item(135, 14).Value = -33473  # N135: <absent> -> -33473

$ws.Cells.Item(136, 8).Value = 100001550  # H136: 83334780 -> 100001550
$ws.Cells.Item(136, 9).Value = 166667500  # I136: 125000860 -> 166667500
$ws.Cells.Item(136, 11).Value = 500002500  # K136: 375002580 -> 500002500
$ws.Cells.Item(136, 13).Value = -499999950  # M136: -375000030 -> -499999950

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 66667784  # H99: 62501096 -> 66667784
$ws.Cells.Item(99, 10).Value = 1274.8334  # J99: 1207 -> 1274.8334
$ws.Cells.Item(99, 12).Value = 1274.8334  # L99: 1207 -> 1274.8334
$ws.Cells.Item(99, 14).Value = -4270.8334  # N99: -4203 -> -4270.8334

$ws.Cells.Item(105, 8).Value = 252473220  # H105: 126236910 -> 252473220
$ws.Cells.Item(105, 9).Value = 336630300  # I105: 168315440 -> 336630300
$ws.Cells.Item(105, 10).Value = 2000  # J105: 1300 -> 2000
$ws.Cells.Item(105, 11).Value = 336630300  # K105: 168315440 -> 336630300
$ws.Cells.Item(105, 12).Value = 2000  # L105: 1300 -> 2000
$ws.Cells.Item(105, 13).Value = -336628553  # M105: -168313693 -> -336628553
$ws.Cells.Item(105, 14).Value = -5494  # N105: -4794 -> -5494

$ws.Cells.Item(114, 8).Value = 0  # H114: 34980 -> 0
$ws.Cells.Item(114, 10).Value = 0  # J114: 34980 -> 0
$ws.Cells.Item(114, 12).Value = 0  # L114: 34980 -> 0
$ws.Cells.Item(114, 14).ClearContents()  # N114: remove (was -43658)

$ws.Cells.Item(137, 8).Value = 37559.832  # H137: 37643.168 -> 37559.832
$ws.Cells.Item(137, 10).Value = 37559.832  # J137: 37643.168 -> 37559.832
$ws.Cells.Item(137, 12).Value = 37559.832  # L137: 37643.168 -> 37559.832
$ws.Cells.Item(137, 14).Value = -47759.832  # N137: -47843.168 -> -47759.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 8531.267  # H58: 7621.706 -> 8531.267
$ws.Cells.Item(58, 9).Value = 1162  # I58: 1081.7142 -> 1162
$ws.Cells.Item(58, 10).Value = 13444.111  # J58: 12199.7 -> 13444.111
$ws.Cells.Item(58, 11).Value = 1162  # K58: 1081.7142 -> 1162
$ws.Cells.Item(58, 12).Value = 13444.111  # L58: 12199.7 -> 13444.111
$ws.Cells.Item(58, 13).Value = -959  # M58: -878.7141999999999 -> -959
$ws.Cells.Item(58, 14).Value = -13850.111  # N58: -12605.7 -> -13850.111

$ws.Cells.Item(107, 8).Value = 1097.3529  # H107: 1295.1538 -> 1097.3529
$ws.Cells.Item(107, 9).Value = 705.8333  # I107: 829.75 -> 705.8333
$ws.Cells.Item(107, 10).Value = 2037  # J107: 2039.8 -> 2037
$ws.Cells.Item(107, 11).Value = 705.8333  # K107: 829.75 -> 705.8333
$ws.Cells.Item(107, 12).Value = 2037  # L107: 2039.8 -> 2037
$ws.Cells.Item(107, 13).Value = 1214.1667  # M107: 1090.25 -> 1214.1667
$ws.Cells.Item(107, 14).Value = -5877  # N107: -5879.8 -> -5877

$ws.Cells.Item(136, 8).Value = 8531.267  # H136: 7621.706 -> 8531.267
$ws.Cells.Item(136, 9).Value = 1162  # I136: 1081.7142 -> 1162
$ws.Cells.Item(136, 10).Value = 13444.111  # J136: 12199.7 -> 13444.111
$ws.Cells.Item(136, 11).Value = 3486  # K136: 3245.1426 -> 3486
$ws.Cells.Item(136, 12).Value = 40332.333  # L136: 36599.10000000001 -> 40332.333
$ws.Cells.Item(136, 13).Value = -936  # M136: -695.1425999999997 -> -936
$ws.Cells.Item(136, 14).Value = -45432.333  # N136: -41699.10000000001 -> -45432.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 940.55554  # H68: 986.5 -> 940.55554
$ws.Cells.Item(68, 10).Value = 738.5714  # J68: 821.25 -> 738.5714
$ws.Cells.Item(68, 12).Value = 2215.7142  # L68: 2463.75 -> 2215.7142
$ws.Cells.Item(68, 14).Value = -3837.7142  # N68: -4085.75 -> -3837.7142

$ws.Cells.Item(71, 8).Value = 940.55554  # H71: 986.5 -> 940.55554
$ws.Cells.Item(71, 10).Value = 738.5714  # J71: 821.25 -> 738.5714
$ws.Cells.Item(71, 12).Value = 6647.1426  # L71: 7391.25 -> 6647.1426
$ws.Cells.Item(71, 14).Value = -14759.1426  # N71: -15503.25 -> -14759.1426

$ws.Cells.Item(107, 8).Value = 8948.166999999999  # H107: 7734.357 -> 8948.166999999999
$ws.Cells.Item(107, 9).Value = 595  # I107: 556.6 -> 595
$ws.Cells.Item(107, 10).Value = 13124.75  # J107: 11722 -> 13124.75
$ws.Cells.Item(107, 11).Value = 1785  # K107: 1669.8 -> 1785
$ws.Cells.Item(107, 12).Value = 39374.25  # L107: 35166 -> 39374.25
$ws.Cells.Item(107, 13).Value = 135  # M107: 250.1999999999998 -> 135
$ws.Cells.Item(107, 14).Value = -43214.25  # N107: -39006 -> -43214.25

$ws.Cells.Item(113, 8).Value = 703.2653  # H113: 695 -> 703.2653
$ws.Cells.Item(113, 9).Value = 641.8  # I113: 635.0476 -> 641.8
$ws.Cells.Item(113, 10).Value = 745.65515  # J113: 735.6129 -> 745.65515
$ws.Cells.Item(113, 11).Value = 1925.4  # K113: 1905.1428 -> 1925.4
$ws.Cells.Item(113, 12).Value = 2236.96545  # L113: 2206.8387 -> 2236.96545
$ws.Cells.Item(113, 13).Value = 244.6000000000001  # M113: 264.8571999999999 -> 244.6000000000001
$ws.Cells.Item(113, 14).Value = -6576.96545  # N113: -6546.8387 -> -6576.96545

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 150001660  # H70: 30003702 -> 150001660
$ws.Cells.Item(70, 9).Value = 125002504  # I70: 22731004 -> 125002504
$ws.Cells.Item(70, 10).Value = 200000000  # J70: 50003620 -> 200000000
$ws.Cells.Item(70, 11).Value = 125002504  # K70: 22731004 -> 125002504
$ws.Cells.Item(70, 12).Value = 200000000  # L70: 50003620 -> 200000000
$ws.Cells.Item(70, 13).Value = -125002234  # M70: -22730734 -> -125002234
$ws.Cells.Item(70, 14).Value = -200000540  # N70: -50004160 -> -200000540

$ws.Cells.Item(73, 8).Value = 150001660  # H73: 30003702 -> 150001660
$ws.Cells.Item(73, 9).Value = 125002504  # I73: 22731004 -> 125002504
$ws.Cells.Item(73, 10).Value = 200000000  # J73: 50003620 -> 200000000
$ws.Cells.Item(73, 11).Value = 125002504  # K73: 22731004 -> 125002504
$ws.Cells.Item(73, 12).Value = 200000000  # L73: 50003620 -> 200000000
$ws.Cells.Item(73, 13).Value = -125001568  # M73: -22730068 -> -125001568
$ws.Cells.Item(73, 14).Value = -200001872  # N73: -50005492 -> -200001872

$ws.Cells.Item(80, 8).Value = 3024.9285  # H80: 3110.6428 -> 3024.9285
$ws.Cells.Item(80, 10).Value = 3705.4443  # J80: 3634.9 -> 3705.4443
$ws.Cells.Item(80, 12).Value = 3705.4443  # L80: 3634.9 -> 3705.4443
$ws.Cells.Item(80, 14).Value = -5701.4443  # N80: -5630.9 -> -5701.4443

$ws.Cells.Item(83, 8).Value = 3024.9285  # H83: 3110.6428 -> 3024.9285
$ws.Cells.Item(83, 10).Value = 3705.4443  # J83: 3634.9 -> 3705.4443
$ws.Cells.Item(83, 12).Value = 18527.2215  # L83: 18174.5 -> 18527.2215
$ws.Cells.Item(83, 14).Value = -28511.2215  # N83: -28158.5 -> -28511.2215

$ws.Cells.Item(122, 8).Value = 1349.5  # H122: 1669 -> 1349.5
$ws.Cells.Item(122, 9).Value = 999  # I122: 1007 -> 999
$ws.Cells.Item(122, 10).Value = 1700  # J122: 2000 -> 1700
$ws.Cells.Item(122, 11).Value = 2997  # K122: 3021 -> 2997
$ws.Cells.Item(122, 12).Value = 5100  # L122: 6000 -> 5100
$ws.Cells.Item(122, 13).Value = -547  # M122: -571 -> -547
$ws.Cells.Item(122, 14).Value = -10000  # N122: -10900 -> -10000

$ws.Cells.Item(132, 8).Value = 7208.625  # H132: 7478.5654 -> 7208.625
$ws.Cells.Item(132, 9).Value = 10253.462  # I132: 11024.583 -> 10253.462
$ws.Cells.Item(132, 11).Value = 30760.386  # K132: 33073.749 -> 30760.386
$ws.Cells.Item(132, 13).Value = -28230.386  # M132: -30543.749 -> -28230.386

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 20240930  # H122: 13891471 -> 20240930
$ws.Cells.Item(122, 9).Value = 28335854  # I122: 15627467 -> 28335854
$ws.Cells.Item(122, 10).Value = 3623.75  # J122: 3500 -> 3623.75
$ws.Cells.Item(122, 11).Value = 85007562  # K122: 46882401 -> 85007562
$ws.Cells.Item(122, 12).Value = 10871.25  # L122: 10500 -> 10871.25
$ws.Cells.Item(122, 13).Value = -85005112  # M122: -46879951 -> -85005112
$ws.Cells.Item(122, 14).Value = -15771.25  # N122: -15400 -> -15771.25

$ws.Cells.Item(135, 8).Value = 35346.668  # H135: 36163.332 -> 35346.668
$ws.Cells.Item(135, 10).Value = 35346.668  # J135: 36163.332 -> 35346.668
$ws.Cells.Item(135, 12).Value = 35346.668  # L135: 36163.332 -> 35346.668
$ws.Cells.Item(135, 14).Value = -45486.668  # N135: -46303.332 -> -45486.668

$ws.Cells.Item(136, 8).Value = 1557.3182  # H136: 1598.1428 -> 1557.3182
$ws.Cells.Item(136, 9).Value = 1409.4445  # I136: 1451.1765 -> 1409.4445
$ws.Cells.Item(136, 11).Value = 4228.333500000001  # K136: 4353.529500000001 -> 4228.333500000001
$ws.Cells.Item(136, 13).Value = -1678.333500000001  # M136: -1803.529500000001 -> -1678.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1850.1613  # H81: 1865.7 -> 1850.1613
$ws.Cells.Item(81, 9).Value = 1478.75  # I81: 1543.6666 -> 1478.75
$ws.Cells.Item(81, 10).Value = 1905.1852  # J81: 1901.4814 -> 1905.1852
$ws.Cells.Item(81, 11).Value = 2957.5  # K81: 3087.3332 -> 2957.5
$ws.Cells.Item(81, 12).Value = 3810.3704  # L81: 3802.9628 -> 3810.3704
$ws.Cells.Item(81, 13).Value = -1896.5  # M81: -2026.3332 -> -1896.5
$ws.Cells.Item(81, 14).Value = -5932.3704  # N81: -5924.962799999999 -> -5932.3704

$ws.Cells.Item(84, 8).Value = 1850.1613  # H84: 1865.7 -> 1850.1613
$ws.Cells.Item(84, 9).Value = 1478.75  # I84: 1543.6666 -> 1478.75
$ws.Cells.Item(84, 10).Value = 1905.1852  # J84: 1901.4814 -> 1905.1852
$ws.Cells.Item(84, 11).Value = 14787.5  # K84: 15436.666 -> 14787.5
$ws.Cells.Item(84, 12).Value = 19051.852  # L84: 19014.814 -> 19051.852
$ws.Cells.Item(84, 13).Value = -9483.5  # M84: -10132.666 -> -9483.5
$ws.Cells.Item(84, 14).Value = -29659.852  # N84: -29622.814 -> -29659.852

$ws.Cells.Item(122, 8).Value = 10418698  # H122: 11365775 -> 10418698
$ws.Cells.Item(122, 9).Value = 14707878  # I122: 16668816 -> 14707878
$ws.Cells.Item(122, 11).Value = 44123634  # K122: 50006448 -> 44123634
$ws.Cells.Item(122, 13).Value = -44121184  # M122: -50003998 -> -44121184

$ws.Cells.Item(135, 8).Value = 64963  # H135: 72495 -> 64963
$ws.Cells.Item(135, 10).Value = 64963  # J135: 72495 -> 64963
$ws.Cells.Item(135, 12).Value = 64963  # L135: 72495 -> 64963
$ws.Cells.Item(135, 14).Value = -75103  # N135: -82635 -> -75103

$ws.Cells.Item(136, 8).Value = 1107.6842  # H136: 1528.5 -> 1107.6842
$ws.Cells.Item(136, 9).Value = 874.4  # I136: 1160.625 -> 874.4
$ws.Cells.Item(136, 10).Value = 1982.5  # J136: 3000 -> 1982.5
$ws.Cells.Item(136, 11).Value = 2623.2  # K136: 3481.875 -> 2623.2
$ws.Cells.Item(136, 12).Value = 5947.5  # L136: 9000 -> 5947.5
$ws.Cells.Item(136, 13).Value = -73.19999999999982  # M136: -931.875 -> -73.19999999999982
$ws.Cells.Item(136, 14).Value = -11047.5  # N136: -14100 -> -11047.5
